$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Row 1: turn the old duplicate-data header row into real column headers ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"

# --- Row 2 / Row 3: fix the swapped register_date / name values ---
$ws.Range("E2").Value = "曰本馬自達MPV"
$ws.Range("B3").Value = "88年04月13曰"

# --- New columns H:N (property_category, category, date, legislator_name, legislator_id, source_file, index) ---
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2011-12-28"
$ws.Range("K2").Value = "王進士"
$ws.Range("L2").Value = 1701
$ws.Range("M2").Value = "tmpf41"
$ws.Range("N2").Value = 47

$ws.Range("H3").Value = "land"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "2011-12-28"
$ws.Range("K3").Value = "王進士"
$ws.Range("L3").Value = 1701
$ws.Range("M3").Value = "tmpf41"
$ws.Range("N3").Value = 48

# --- Match formatting: header row (H1:N1) picks up the bold/bordered style used by B1:G1,
#     data rows (H2:N2 / H3:N3) pick up the plain style used by B2:G2 / B3:G3 ---
$ws.Range("G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

$ws.Range("G2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)

$ws.Range("G3").Copy()
$ws.Range("H3:N3").PasteSpecial(-4122)

$excel.CutCopyMode = $false
